$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column H: copy the date-header format from G2, then write values ---
$ws.Range("G2").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("H2").Value = 43995
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 3
$ws.Range("H10").Value = 4
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 4
$ws.Range("H13").Value = 6
$ws.Range("H14").Value = 8
$ws.Range("H15").Value = 1
$ws.Range("H16").Value = 5
$ws.Range("H17").Value = 4
$ws.Range("H18").Value = 10
$ws.Range("H19").Value = 9
$ws.Range("H20").Value = 6
$ws.Range("H21").Value = 7
$ws.Range("H22").Value = 3
$ws.Range("H23").Value = 9
$ws.Range("H24").Value = 4
$ws.Range("H25").Value = 8
$ws.Range("H26").Value = 6
$ws.Range("H27").Value = 5
$ws.Range("H28").Value = 6
$ws.Range("H29").Value = 5
$ws.Range("H30").Value = 10
$ws.Range("H31").Value = 8
$ws.Range("H32").Value = 5
$ws.Range("H33").Value = 7
$ws.Range("H34").Value = 6
$ws.Range("H35").Value = 7
$ws.Range("H36").Value = 9
$ws.Range("H37").Value = 11
$ws.Range("H38").Value = 7
$ws.Range("H39").Value = 6
$ws.Range("H40").Value = 8
$ws.Range("H41").Value = 10
$ws.Range("H42").Value = 7
$ws.Range("H43").Value = 19
$ws.Range("H44").Value = 12
$ws.Range("H45").Value = 11
$ws.Range("H46").Value = 14
$ws.Range("H47").Value = 29
$ws.Range("H48").Value = 16
$ws.Range("H49").Value = 22
$ws.Range("H50").Value = 30
$ws.Range("H51").Value = 19
$ws.Range("H52").Value = 25
$ws.Range("H53").Value = 38
$ws.Range("H54").Value = 35
$ws.Range("H55").Value = 32
$ws.Range("H56").Value = 48
$ws.Range("H57").Value = 45
$ws.Range("H58").Value = 60
$ws.Range("H59").Value = 60
$ws.Range("H60").Value = 65
$ws.Range("H61").Value = 60
$ws.Range("H62").Value = 69
$ws.Range("H63").Value = 85
$ws.Range("H64").Value = 94
$ws.Range("H65").Value = 85
$ws.Range("H66").Value = 102
$ws.Range("H67").Value = 113
$ws.Range("H68").Value = 110
$ws.Range("H69").Value = 119
$ws.Range("H70").Value = 121
$ws.Range("H71").Value = 108
$ws.Range("H72").Value = 133
$ws.Range("H73").Value = 129
$ws.Range("H74").Value = 135
$ws.Range("H75").Value = 132
$ws.Range("H76").Value = 99
$ws.Range("H77").Value = 118
$ws.Range("H78").Value = 121
$ws.Range("H79").Value = 114
$ws.Range("H80").Value = 127
$ws.Range("H81").Value = 97
$ws.Range("H82").Value = 29

# --- H1 grand-total formula ---
$ws.Range("H1").Formula = "=SUM(H3:H82)"

$excel.Calculate()

# --- Restore view / selection state ---
$ws.Range("M13").Select()
